# Insert a new slide (Comparison of garlic vs non-garlic food combos)
# right before the final "Where do we go from here?..." slide.

$p = $ppt.ActivePresentation

# Add the new slide in position 7 (pushing the current slide 7 to slide 8).
# Use the "Title Only" layout so PowerPoint materializes a real Title
# placeholder for us, then switch the slide over to the "Blank" layout
# (the layout used by the rest of the deck's content slides) while keeping
# the title shape that's already on the slide.
$newSlide = $p.Slides.Add(7, 6)

$blankLayout = $p.SlideMaster.CustomLayouts.Item(7)
$newSlide.CustomLayout = $blankLayout

# Title text
$newSlide.Shapes.Title.TextFrame.TextRange.Text = "Comparison of garlic vs non-garlic food combos"

# EMU -> point helper (PowerPoint COM shape geometry is expressed in points)
$emuPerPt = 12700

$tblLeft = 1608488 / $emuPerPt
$tblTop = 2355961 / $emuPerPt
$tblWidth = 7940951 / $emuPerPt
$tblHeight = 3337560 / $emuPerPt

$tblShape = $newSlide.Shapes.AddTable(9, 4, $tblLeft, $tblTop, $tblWidth, $tblHeight)
$tblShape.Name = "Table 3"

$tbl = $tblShape.Table

# Match the author's hand-picked column widths
$colWidths = @(3280200, 978869, 1857626, 1824256)
for ($c = 1; $c -le 4; $c++) {
    $tbl.Columns.Item($c).Width = $colWidths[$c - 1] / $emuPerPt
}

# Apply the "Medium Style 2 - Accent 2" table style used in the authored deck
$tbl.ApplyStyle("{5FD0F851-EC5A-4D38-B0AD-8093EC10F338}")

# Header row
$tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Food Combo"
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Count"
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "Average_eat"
$tbl.Cell(1, 4).Shape.TextFrame.TextRange.Text = "Eat_target_pct"

# Data rows
$rows = @(
    @("Mackerel_Garlic", "108", "19.4", "0.32"),
    @("Mackerel", "36", "18.2", "0.30"),
    @("Herring_Garlic", "33", "12.2", "0.36"),
    @("Herring", "89", "11.4", "0.48"),
    @("Saury Garlic", "103", "27.1", "0.41"),
    @("Saury", "16", "24.1", "0.35"),
    @("Blue Runner_Squid Garlic", "67", "28.3", "0.39"),
    @("Blue Runner_squid", "10", "20.6", "0.30")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Text = $row[0]
    $tbl.Cell($r, 2).Shape.TextFrame.TextRange.Text = $row[1]
    $tbl.Cell($r, 3).Shape.TextFrame.TextRange.Text = $row[2]
    $tbl.Cell($r, 4).Shape.TextFrame.TextRange.Text = $row[3]
}
